# Insert a new "Tag Line" slide right after the title slide ("Learning quest"),
# pushing all of the existing content slides down by one position.
#
# This mirrors using PowerPoint's "New Slide" command while the first slide is
# selected: a new slide is created at index 2 using the same "Title and
# Content" layout already used by the rest of the deck, and then its title
# and body placeholders are filled in.

$p = $ppt.ActivePresentation

# "Title and Content" is CustomLayouts.Item(2) on this deck's slide master
# (it's the layout already used by every non-title slide).
$layout = $p.SlideMaster.CustomLayouts.Item(2)

$newSlide = $p.Slides.AddSlide(2, $layout)

# Title placeholder.
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Tag Line`t`t`t`t"

# Body / content placeholder.
$body = $newSlide.Shapes.Item(2)
$bodyText = "Learning quest is a colourful playground adventure where even the sky isn’t the limit. Run and jump from swings to clouds picking up letters to spell words and progress through levels.`r`r`r`r`r`r"
$body.TextFrame.TextRange.Text = $bodyText
